# Add two new rows of paper data (rows 11-12) that were previously blank,
# add the two new hyperlinks that go with them, give the new arXiv link
# cell the "followed hyperlink" (purple, underlined) look, widen column A,
# and move the sheet's viewport/selection the way the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 11: CaEGCN paper -------------------------------------------------
$ws.Range("A11").Value = "CaEGCN: Cross-Attention Fusion based Enhanced Graph Convolutional Network forClustering"
$ws.Range("B11").Value = "IEEE Transactions on Knowledge and Data Engineering"
$ws.Range("C11").Value = "交叉注意卷积网络增强"
$ws.Range("D11").Value = "是"
$ws.Range("E11").Value = "huogy/CaEGCN (github.com)"
$ws.Hyperlinks.Add($ws.Range("E11"), "https://github.com/huogy/CaEGCN", "", "https://github.com/huogy/CaEGCN", "huogy/CaEGCN (github.com)")

# --- Row 12: Rethinking Cross-Attention paper ------------------------------
$ws.Range("A12").Value = "Rethinking Cross-Attention for Infrared and Visible Image Fusion"
$ws.Range("B12").Value = "https://arxiv.org/abs/2401.11675"
$ws.Hyperlinks.Add($ws.Range("B12"), "https://arxiv.org/abs/2401.11675", "", "", "https://arxiv.org/abs/2401.11675")
# This link was already visited by the author, so it uses the purple
# "FollowedHyperlink" colour instead of the usual blue one.
$ws.Range("B12").Font.Color = 8388736
$ws.Range("B12").Font.Underline = 2
$ws.Range("C12").Value = "注意机制"
$ws.Range("D12").Value = "是"
$ws.Range("E12").Value = "https://github.com/Linfeng-Tang/PSFusion"

# --- Column A is now wider to fit the longer paper titles -----------------
$ws.Columns.Item(1).ColumnWidth = 89.4

# --- Viewport: no more frozen "topLeftCell", selection moved to B10 -------
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("B10").Select()
